# Update the start/end time values on Sheet1. These were recorded via
# a timestamp-producing routine (e.g. Get-Date / epoch-ish ticks) and are
# being refreshed with a new pair of readings; the dependent formulas in
# B3 (difference) and B4 (difference in seconds) recalculate on their own.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$startTime = [double]"1.5018674201649999E+18"
$endTime   = [double]"1.5018674291219999E+18"

$ws.Range("B1").Value = $startTime
$ws.Range("B2").Value = $endTime

# The cursor ends up resting on the "difference" cell after the refresh.
$ws.Range("B3").Select() | Out-Null
